# Auto-generated edit script: updates Leve profit-calculation values
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to reflect refreshed
# market-board price data pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (26 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 61388.168
$ws.Range("J87").Value = 61388.168
$ws.Range("L87").Value = 61388.168
$ws.Range("N87").Value = -63884.168
$ws.Range("H90").Value = 61388.168
$ws.Range("J90").Value = 61388.168
$ws.Range("L90").Value = 184164.504
$ws.Range("N90").Value = -196644.504
$ws.Range("H132").Value = 1649.5
$ws.Range("I132").Value = 1649.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4948.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2418.5
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 4899.3706
$ws.Range("I137").Value = 1771.9474
$ws.Range("K137").Value = 5315.8422
$ws.Range("M137").Value = -2765.8422
$ws.Range("H138").Value = 2246.64
$ws.Range("I138").Value = 1911.1333
$ws.Range("J138").Value = 2749.9
$ws.Range("K138").Value = 5733.3999
$ws.Range("L138").Value = 8249.700000000001
$ws.Range("M138").Value = -593.3999000000003
$ws.Range("N138").Value = -18529.7

# --- Sheet: ARM (38 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2302705.5
$ws.Range("I2").Value = 2423852
$ws.Range("K2").Value = 2423852
$ws.Range("M2").Value = -2423739
$ws.Range("H32").Value = 25361.195
$ws.Range("I32").Value = 25171.877
$ws.Range("K32").Value = 25171.877
$ws.Range("M32").Value = -24884.877
$ws.Range("H97").Value = 1160810.1
$ws.Range("I97").Value = 1613762
$ws.Range("J97").Value = 3266.3333
$ws.Range("K97").Value = 1613762
$ws.Range("L97").Value = 3266.3333
$ws.Range("M97").Value = -1613266
$ws.Range("N97").Value = -4258.3333
$ws.Range("H101").Value = 54285.285
$ws.Range("J101").Value = 54285.285
$ws.Range("L101").Value = 54285.285
$ws.Range("N101").Value = -60775.285
$ws.Range("H102").Value = 22821472
$ws.Range("I102").Value = 3740.2856
$ws.Range("J102").Value = 62752504
$ws.Range("K102").Value = 3740.2856
$ws.Range("L102").Value = 62752504
$ws.Range("M102").Value = -2118.2856
$ws.Range("N102").Value = -62755748
$ws.Range("H116").Value = 2302705.5
$ws.Range("I116").Value = 2423852
$ws.Range("K116").Value = 2423852
$ws.Range("M116").Value = -2421558
$ws.Range("H132").Value = 3403.3684
$ws.Range("I132").Value = 3026.4443
$ws.Range("K132").Value = 9079.332900000001
$ws.Range("M132").Value = -6549.332900000001
$ws.Range("H133").Value = 117499.5
$ws.Range("J133").Value = 117499.5
$ws.Range("L133").Value = 117499.5
$ws.Range("N133").Value = -122559.5

# --- Sheet: BSM (16 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2302705.5
$ws.Range("I3").Value = 2423852
$ws.Range("K3").Value = 2423852
$ws.Range("M3").Value = -2423738
$ws.Range("H105").Value = 62520588
$ws.Range("I105").Value = 90935624
$ws.Range("K105").Value = 90935624
$ws.Range("M105").Value = -90933877
$ws.Range("H107").Value = 1194.5938
$ws.Range("I107").Value = 1141.0454
$ws.Range("K107").Value = 1141.0454
$ws.Range("M107").Value = 778.9546
$ws.Range("H134").Value = 3813.0613
$ws.Range("J134").Value = 6929.077
$ws.Range("L134").Value = 20787.231
$ws.Range("N134").Value = -25857.231

# --- Sheet: CRP (16 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H103").Value = 44605.25
$ws.Range("J103").Value = 74899.5
$ws.Range("L103").Value = 74899.5
$ws.Range("N103").Value = -77243.5
$ws.Range("H112").Value = 93734.664
$ws.Range("J112").Value = 93734.664
$ws.Range("L112").Value = 93734.664
$ws.Range("N112").Value = -96688.664
$ws.Range("H132").Value = 2358.7307
$ws.Range("I132").Value = 1166.3914
$ws.Range("K132").Value = 3499.1742
$ws.Range("M132").Value = -969.1741999999999
$ws.Range("H134").Value = 1965.5
$ws.Range("I134").Value = 1657.3334
$ws.Range("K134").Value = 4972.0002
$ws.Range("M134").Value = -2437.0002

# --- Sheet: CUL (42 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 3759.8
$ws.Range("I93").Value = 1399.5
$ws.Range("J93").Value = 5333.3335
$ws.Range("K93").Value = 4198.5
$ws.Range("L93").Value = 16000.0005
$ws.Range("M93").Value = -2326.5
$ws.Range("N93").Value = -19744.0005
$ws.Range("H94").Value = 4960
$ws.Range("J94").Value = 4960
$ws.Range("L94").Value = 14880
$ws.Range("N94").Value = -16232
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H100").Value = 4499.5
$ws.Range("J100").Value = 7999
$ws.Range("L100").Value = 23997
$ws.Range("N100").Value = -25619
$ws.Range("H101").Value = 7999
$ws.Range("I101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("M104").ClearContents()
$ws.Range("N104").ClearContents()
$ws.Range("H120").Value = 859.6667
$ws.Range("I120").Value = 859.6667
$ws.Range("K120").Value = 2579.0001
$ws.Range("M120").Value = 2258.9999
$ws.Range("H131").Value = 10106206
$ws.Range("J131").Value = 8169.65
$ws.Range("L131").Value = 24508.95
$ws.Range("N131").Value = -34588.95
$ws.Range("H140").Value = 1339.9
$ws.Range("J140").Value = 3033
$ws.Range("L140").Value = 9099
$ws.Range("N140").Value = -19459

# --- Sheet: GSM (29 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6680
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 6680
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 6680
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -7220
$ws.Range("H73").Value = 6680
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 6680
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 6680
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -8552
$ws.Range("H80").Value = 5973
$ws.Range("I80").Value = 4999
$ws.Range("K80").Value = 4999
$ws.Range("M80").Value = -4001
$ws.Range("H83").Value = 5973
$ws.Range("I83").Value = 4999
$ws.Range("K83").Value = 24995
$ws.Range("M83").Value = -20003
$ws.Range("H132").Value = 5844.0303
$ws.Range("I132").Value = 4465.2915
$ws.Range("J132").Value = 9520.666999999999
$ws.Range("K132").Value = 13395.8745
$ws.Range("L132").Value = 28562.001
$ws.Range("M132").Value = -10865.8745
$ws.Range("N132").Value = -33622.001

# --- Sheet: LTW (64 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 833.75
$ws.Range("I16").Value = 856
$ws.Range("K16").Value = 856
$ws.Range("M16").Value = -686
$ws.Range("H46").Value = 7771
$ws.Range("I46").Value = 1500.5
$ws.Range("K46").Value = 1500.5
$ws.Range("M46").Value = -1312.5
$ws.Range("H68").Value = 2188.3333
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 2565
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 2565
$ws.Range("M68").Value = -1251
$ws.Range("N68").Value = -4063
$ws.Range("H71").Value = 2188.3333
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 2565
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 12825
$ws.Range("M71").Value = -6256
$ws.Range("N71").Value = -20313
$ws.Range("H82").Value = 1174.9375
$ws.Range("I82").Value = 1766.6666
$ws.Range("J82").Value = 819.9
$ws.Range("K82").Value = 1766.6666
$ws.Range("L82").Value = 819.9
$ws.Range("M82").Value = -1405.6666
$ws.Range("N82").Value = -1541.9
$ws.Range("H85").Value = 1174.9375
$ws.Range("I85").Value = 1766.6666
$ws.Range("J85").Value = 819.9
$ws.Range("K85").Value = 1766.6666
$ws.Range("L85").Value = 819.9
$ws.Range("M85").Value = -518.6666
$ws.Range("N85").Value = -3315.9
$ws.Range("H93").Value = 3081.625
$ws.Range("I93").Value = 2762.9412
$ws.Range("J93").Value = 3855.5715
$ws.Range("K93").Value = 2762.9412
$ws.Range("L93").Value = 3855.5715
$ws.Range("M93").Value = -1514.9412
$ws.Range("N93").Value = -6351.5715
$ws.Range("H100").Value = 10872690
$ws.Range("I100").Value = 83334664
$ws.Range("J100").Value = 3393.05
$ws.Range("K100").Value = 83334664
$ws.Range("L100").Value = 3393.05
$ws.Range("M100").Value = -83334123
$ws.Range("N100").Value = -4475.05
$ws.Range("H132").Value = 6697.8
$ws.Range("I132").Value = 2996.3333
$ws.Range("J132").Value = 12250
$ws.Range("K132").Value = 8988.999899999999
$ws.Range("L132").Value = 36750
$ws.Range("M132").Value = -6458.999899999999
$ws.Range("N132").Value = -41810
$ws.Range("H136").Value = 4954.3335
$ws.Range("I136").Value = 3622.7334
$ws.Range("J136").Value = 8283.333000000001
$ws.Range("K136").Value = 10868.2002
$ws.Range("L136").Value = 24849.999
$ws.Range("M136").Value = -8318.200199999999
$ws.Range("N136").Value = -29949.999

# --- Sheet: WVR (26 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 30749.5
$ws.Range("J39").Value = 39899.2
$ws.Range("L39").Value = 39899.2
$ws.Range("N39").Value = -40725.2
$ws.Range("H43").Value = 64072.215
$ws.Range("J43").Value = 64728.273
$ws.Range("L43").Value = 64728.273
$ws.Range("N43").Value = -65026.273
$ws.Range("H51").Value = 538126.7
$ws.Range("I51").Value = 877667.7
$ws.Range("K51").Value = 877667.7
$ws.Range("M51").Value = -877157.7
$ws.Range("H126").Value = 3652.121
$ws.Range("I126").Value = 3463.9583
$ws.Range("J126").Value = 4153.8887
$ws.Range("K126").Value = 10391.8749
$ws.Range("L126").Value = 12461.6661
$ws.Range("M126").Value = -7921.874899999999
$ws.Range("N126").Value = -17401.6661
$ws.Range("H136").Value = 3022.0557
$ws.Range("I136").Value = 1365.4736
$ws.Range("J136").Value = 4873.5293
$ws.Range("K136").Value = 4096.4208
$ws.Range("L136").Value = 14620.5879
$ws.Range("M136").Value = -1546.4208
$ws.Range("N136").Value = -19720.5879

